$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds locale-formatted numbers-as-text (e.g. "41.983.35",
# "6.80"); force text formatting before assignment so Excel does not
# auto-coerce/round them into numeric values.

# Row 36 and 37: coin identities swap (Celestia <-> WEMIXToken), with new price/volume data
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.93"
$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  -0.02%  "

# Remaining price / volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.983.35"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.254.50"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.63"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.82"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.86"
$ws.Range("E10").Value = "  +5.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.80"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.601.99"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.39"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.248.23"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.839.96"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.94"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.35"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.69"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.49"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.43"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.40"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.28"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0733"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.958.28"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.55"
$ws.Range("E45").Value = "  -5.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.91"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.83"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.23"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.92"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -1.91%  "
